$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap A2 <-> B2
$a2 = $ws.Range("A2").Value2
$b2 = $ws.Range("B2").Value2
$ws.Range("A2").Value = $b2
$ws.Range("B2").Value = $a2

# Swap C2 <-> D2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$ws.Range("C2").Value = $d2
$ws.Range("D2").Value = $c2

# Move selection to D3 (as if the user pressed Enter after editing D2)
$ws.Range("D3").Select()
